# Auto-generated Excel COM-interop edit script
# Applies numeric cell updates (cached literal values) per the target diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(111, 8).Value = 1177.5714
$ws.Cells.Item(111, 10).Value = 1788.6666
$ws.Cells.Item(111, 12).Value = 5365.9998
$ws.Cells.Item(111, 14).Value = -11499.9998

$ws.Cells.Item(135, 8).Value = 1930.4736
$ws.Cells.Item(135, 9).Value = 1911.9333
$ws.Cells.Item(135, 10).Value = 2000
$ws.Cells.Item(135, 11).Value = 17207.3997
$ws.Cells.Item(135, 12).Value = 18000
$ws.Cells.Item(135, 13).Value = -14672.3997
$ws.Cells.Item(135, 14).Value = -23070

$ws.Cells.Item(137, 8).Value = 1201.1063
$ws.Cells.Item(137, 9).Value = 787.1053000000001
$ws.Cells.Item(137, 10).Value = 1482.0358
$ws.Cells.Item(137, 11).Value = 2361.3159
$ws.Cells.Item(137, 12).Value = 4446.107400000001
$ws.Cells.Item(137, 13).Value = 188.6840999999999
$ws.Cells.Item(137, 14).Value = -9546.107400000001

$ws.Cells.Item(138, 8).Value = 1845.18
$ws.Cells.Item(138, 10).Value = 2421.2207
$ws.Cells.Item(138, 12).Value = 7263.6621
$ws.Cells.Item(138, 14).Value = -17543.6621


$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 1177.9756
$ws.Cells.Item(74, 9).Value = 1258.9615
$ws.Cells.Item(74, 10).Value = 1037.6
$ws.Cells.Item(74, 11).Value = 1258.9615
$ws.Cells.Item(74, 12).Value = 1037.6
$ws.Cells.Item(74, 13).Value = -384.9614999999999
$ws.Cells.Item(74, 14).Value = -2785.6

$ws.Cells.Item(77, 8).Value = 1177.9756
$ws.Cells.Item(77, 9).Value = 1258.9615
$ws.Cells.Item(77, 10).Value = 1037.6
$ws.Cells.Item(77, 11).Value = 6294.807499999999
$ws.Cells.Item(77, 12).Value = 5188
$ws.Cells.Item(77, 13).Value = -1926.807499999999
$ws.Cells.Item(77, 14).Value = -13924

$ws.Cells.Item(122, 8).Value = 18066.666
$ws.Cells.Item(122, 9).Value = 18066.666
$ws.Cells.Item(122, 11).Value = 54199.99800000001
$ws.Cells.Item(122, 13).Value = -51749.99800000001

$ws.Cells.Item(133, 8).Value = 50981.637
$ws.Cells.Item(133, 10).Value = 50981.637
$ws.Cells.Item(133, 12).Value = 50981.637
$ws.Cells.Item(133, 14).Value = -56041.637


$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 2094.362
$ws.Cells.Item(134, 9).Value = 1286.2709
$ws.Cells.Item(134, 10).Value = 5973.2
$ws.Cells.Item(134, 11).Value = 3858.8127
$ws.Cells.Item(134, 12).Value = 17919.6
$ws.Cells.Item(134, 13).Value = -1323.8127
$ws.Cells.Item(134, 14).Value = -22989.6


$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(4, 8).Value = 336833340
$ws.Cells.Item(4, 9).Value = 500000
$ws.Cells.Item(4, 11).Value = 500000
$ws.Cells.Item(4, 13).Value = -499888

$ws.Cells.Item(16, 8).Value = 690
$ws.Cells.Item(16, 9).Value = 586.6667
$ws.Cells.Item(16, 10).Value = 1000
$ws.Cells.Item(16, 11).Value = 586.6667
$ws.Cells.Item(16, 12).Value = 1000
$ws.Cells.Item(16, 13).Value = -299.6667
$ws.Cells.Item(16, 14).Value = -1574

$ws.Cells.Item(31, 8).Value = 1197.94
$ws.Cells.Item(31, 9).Value = 883.5
$ws.Cells.Item(31, 10).Value = 1512.38
$ws.Cells.Item(31, 11).Value = 883.5
$ws.Cells.Item(31, 12).Value = 1512.38
$ws.Cells.Item(31, 13).Value = -588.5
$ws.Cells.Item(31, 14).Value = -2102.38

$ws.Cells.Item(34, 8).Value = 1197.94
$ws.Cells.Item(34, 9).Value = 883.5
$ws.Cells.Item(34, 10).Value = 1512.38
$ws.Cells.Item(34, 11).Value = 883.5
$ws.Cells.Item(34, 12).Value = 1512.38
$ws.Cells.Item(34, 13).Value = -681.5
$ws.Cells.Item(34, 14).Value = -1916.38

$ws.Cells.Item(113, 8).Value = 690
$ws.Cells.Item(113, 9).Value = 586.6667
$ws.Cells.Item(113, 10).Value = 1000
$ws.Cells.Item(113, 11).Value = 586.6667
$ws.Cells.Item(113, 12).Value = 1000
$ws.Cells.Item(113, 13).Value = 1583.3333
$ws.Cells.Item(113, 14).Value = -5340

$ws.Cells.Item(122, 8).Value = 803.9091
$ws.Cells.Item(122, 9).Value = 649.1429000000001
$ws.Cells.Item(122, 10).Value = 1074.75
$ws.Cells.Item(122, 11).Value = 1947.4287
$ws.Cells.Item(122, 12).Value = 3224.25
$ws.Cells.Item(122, 13).Value = 502.5712999999998
$ws.Cells.Item(122, 14).Value = -8124.25

$ws.Cells.Item(132, 8).Value = 2534.2856
$ws.Cells.Item(132, 9).Value = 1839.7084
$ws.Cells.Item(132, 10).Value = 4049.7273
$ws.Cells.Item(132, 11).Value = 5519.1252
$ws.Cells.Item(132, 12).Value = 12149.1819
$ws.Cells.Item(132, 13).Value = -2989.1252
$ws.Cells.Item(132, 14).Value = -17209.1819


$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 11485.607
$ws.Cells.Item(4, 9).Value = 199.25
$ws.Cells.Item(4, 10).Value = 13366.667
$ws.Cells.Item(4, 11).Value = 597.75
$ws.Cells.Item(4, 12).Value = 40100.001
$ws.Cells.Item(4, 13).Value = -485.75
$ws.Cells.Item(4, 14).Value = -40324.001

$ws.Cells.Item(39, 8).Value = 8837.031999999999
$ws.Cells.Item(39, 10).Value = 8837.031999999999
$ws.Cells.Item(39, 12).Value = 26511.096
$ws.Cells.Item(39, 14).Value = -27099.096

$ws.Cells.Item(107, 8).Value = 1110.12
$ws.Cells.Item(107, 9).Value = 381.66666
$ws.Cells.Item(107, 10).Value = 1209.4546
$ws.Cells.Item(107, 11).Value = 1144.99998
$ws.Cells.Item(107, 12).Value = 3628.3638
$ws.Cells.Item(107, 13).Value = 775.0000199999999
$ws.Cells.Item(107, 14).Value = -7468.3638

$ws.Cells.Item(131, 8).Value = 2685.1667
$ws.Cells.Item(131, 9).Value = 526
$ws.Cells.Item(131, 10).Value = 3033.4194
$ws.Cells.Item(131, 11).Value = 1578
$ws.Cells.Item(131, 12).Value = 9100.2582
$ws.Cells.Item(131, 13).Value = 3462
$ws.Cells.Item(131, 14).Value = -19180.2582


$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(5, 8).Value = 700
$ws.Cells.Item(5, 9).Value = 400
$ws.Cells.Item(5, 10).Value = 1000
$ws.Cells.Item(5, 11).Value = 400
$ws.Cells.Item(5, 12).Value = 1000
$ws.Cells.Item(5, 13).Value = -288
$ws.Cells.Item(5, 14).Value = -1224

$ws.Cells.Item(102, 8).Value = 4991.5557
$ws.Cells.Item(102, 9).Value = 1602.5
$ws.Cells.Item(102, 10).Value = 7702.8
$ws.Cells.Item(102, 11).Value = 1602.5
$ws.Cells.Item(102, 12).Value = 7702.8
$ws.Cells.Item(102, 13).Value = 19.5
$ws.Cells.Item(102, 14).Value = -10946.8

$ws.Cells.Item(132, 8).Value = 4261.3184
$ws.Cells.Item(132, 9).Value = 3869.375
$ws.Cells.Item(132, 10).Value = 5306.5
$ws.Cells.Item(132, 11).Value = 11608.125
$ws.Cells.Item(132, 12).Value = 15919.5
$ws.Cells.Item(132, 13).Value = -9078.125
$ws.Cells.Item(132, 14).Value = -20979.5


$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(2, 8).Value = 66668.664
$ws.Cells.Item(2, 10).Value = 66668.664
$ws.Cells.Item(2, 12).Value = 66668.664
$ws.Cells.Item(2, 14).Value = -66892.664

$ws.Cells.Item(40, 8).Value = 4397.222
$ws.Cells.Item(40, 9).Value = 2800
$ws.Cells.Item(40, 10).Value = 4491.1763
$ws.Cells.Item(40, 11).Value = 2800
$ws.Cells.Item(40, 12).Value = 4491.1763
$ws.Cells.Item(40, 13).Value = -2664
$ws.Cells.Item(40, 14).Value = -4763.1763

$ws.Cells.Item(100, 8).Value = 2420.6553
$ws.Cells.Item(100, 9).Value = 1818.091
$ws.Cells.Item(100, 10).Value = 2788.889
$ws.Cells.Item(100, 11).Value = 1818.091
$ws.Cells.Item(100, 12).Value = 2788.889
$ws.Cells.Item(100, 13).Value = -1277.091
$ws.Cells.Item(100, 14).Value = -3870.889

$ws.Cells.Item(122, 8).Value = 3874.9375
$ws.Cells.Item(122, 9).Value = 3000
$ws.Cells.Item(122, 10).Value = 3933.2666
$ws.Cells.Item(122, 11).Value = 9000
$ws.Cells.Item(122, 12).Value = 11799.7998
$ws.Cells.Item(122, 13).Value = -6550
$ws.Cells.Item(122, 14).Value = -16699.7998

$ws.Cells.Item(136, 8).Value = 5344.5
$ws.Cells.Item(136, 9).Value = 1929.6471
$ws.Cells.Item(136, 10).Value = 10622
$ws.Cells.Item(136, 11).Value = 5788.9413
$ws.Cells.Item(136, 12).Value = 31866
$ws.Cells.Item(136, 13).Value = -3238.9413
$ws.Cells.Item(136, 14).Value = -36966


$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 101369.4
$ws.Cells.Item(122, 9).Value = 144170.58
$ws.Cells.Item(122, 10).Value = 1500
$ws.Cells.Item(122, 11).Value = 432511.74
$ws.Cells.Item(122, 12).Value = 4500
$ws.Cells.Item(122, 13).Value = -430061.74
$ws.Cells.Item(122, 14).Value = -9400

$ws.Cells.Item(132, 8).Value = 17859408
$ws.Cells.Item(132, 9).Value = 21740808
$ws.Cells.Item(132, 11).Value = 65222424
$ws.Cells.Item(132, 13).Value = -65219894

$ws.Cells.Item(136, 8).Value = 25719996
$ws.Cells.Item(136, 9).Value = 37148944
$ws.Cells.Item(136, 11).Value = 111446832
$ws.Cells.Item(136, 13).Value = -111444282

